$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Sensitivity tests now get their own fuel entry: insert a new row above the
# current "syngas - wood" row (17), which pushes that row (and its comment)
# down to row 18. Re-populate the new row 17 with the "syngas - wood" data
# and relabel the shifted row 18 as the new "syngas - NREL" entry.
# ---------------------------------------------------------------------------

# Grab the comment that currently lives on A17 ("syngas - wood") before we
# disturb anything - we'll need to re-home it on A18 once the shift happens,
# since this engine's row Insert() does not relocate comments on its own.
$oldComment = $ws.Range("A17").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# Insert a new row above row 17; old row 17 (HHV=21.5 etc.) becomes row 18.
$ws.Rows("17:17").Insert()

# New row 17: re-enter the "syngas - wood" fuel data (same figures as before).
$ws.Range("A17").Value = "syngas - wood"
$ws.Range("B17").Value = 21.5
$ws.Range("C17").Value = 21.5
$ws.Range("D17").Value = 1.3
$ws.Range("F17").NumberFormat = "0.0000"
$ws.Range("F17").Formula = "=D17*`$F`$2"

# Row 18 (the shifted original row): rename it to the new "syngas - NREL" fuel.
$ws.Range("A18").Value = "syngas - NREL"

# Re-attach the Swanson 2010 comment on row 18 (where it ends up after save).
$newComment = $ws.Range("A18").AddComment($commentText)

# Row 20 is now "syngas - PNNL" (shifted down from row 19) - it previously had
# no emissions-factor column; give it the same F formula as its neighbours.
$ws.Range("F20").NumberFormat = "0.0000"
$ws.Range("F20").Formula = "=D20*`$F`$2"

# Match the final selection left behind in the saved workbook.
$ws.Range("E26").Select()
